# Ortografia: eliminar signos de apertura ( Spanish inverted ? and ! ) y
# tildes (vocales acentuadas) en todo el texto de todas las hojas del libro.
# La letra "n with tilde" y el resto del contenido (emojis, puntuacion
# normal, numeros, formulas, etc.) se deja intacta; solo se normalizan las
# vocales con tilde/diéresis y se retiran los signos de apertura de
# interrogacion/exclamacion invertidos.

$wb = $excel.ActiveWorkbook

# Pares [caracter a buscar -> reemplazo]. Usamos codigos Unicode (vs. literales)
# para evitar cualquier problema de codificacion al guardar/leer este script.
$pairs = @(
    @([char]0x00E1, "a"),   # U+00E1 a-acute (lowercase)
    @([char]0x00E9, "e"),   # U+00E9 e-acute (lowercase)
    @([char]0x00ED, "i"),   # U+00ED i-acute (lowercase)
    @([char]0x00F3, "o"),   # U+00F3 o-acute (lowercase)
    @([char]0x00FA, "u"),   # U+00FA u-acute (lowercase)
    @([char]0x00C1, "A"),   # U+00C1 A-acute (uppercase)
    @([char]0x00C9, "E"),   # U+00C9 E-acute (uppercase)
    @([char]0x00CD, "I"),   # U+00CD I-acute (uppercase)
    @([char]0x00D3, "O"),   # U+00D3 O-acute (uppercase)
    @([char]0x00DA, "U"),   # U+00DA U-acute (uppercase)
    @([char]0x00BF, ""),    # U+00BF inverted question mark
    @([char]0x00A1, "")     # U+00A1 inverted exclamation mark
)

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    if ($used -eq $null) { continue }

    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2

            if ($val -eq $null) { continue }
            if ($val.GetType().Name -ne "String") { continue }

            $orig = $val
            $new = $val
            foreach ($pair in $pairs) {
                $new = $new.Replace($pair[0], $pair[1])
            }

            if ($new -ne $orig) {
                $cell.Value = $new
            }
        }
    }
}
